$wb = $excel.ActiveWorkbook

# On the "numeric" sheet, fill in the missing id value and move the
# selection to where the user was last working (E3).
$wsNumeric = $wb.Worksheets.Item("numeric")
$wsNumeric.Activate()
$wsNumeric.Range("E2").Value = "*"
$wsNumeric.Range("E3").Select()

# The "string" sheet is the one that should be active/selected when the
# workbook is reopened (previously "drop" was active).
$wsString = $wb.Worksheets.Item("string")
$wsString.Activate()
